$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(3.108436580699082, 1.2902709029944714, -0.3342445738815883, 2.1952681219794457)"
$ws.Range("C2").Value = "NIG(1.1375390676077846, 0.8673294009355148, 3.6796744990945363, 4.878560135865573)"
$ws.Range("D2").Value = "NCT(2.158436442253288, 1.4962217573916945, -0.6980033242432776, 2.192834432140004)"
$ws.Range("E2").Value = "NIG(1.5718075567915843, 1.2653873044936532, 3.7033568041349727, 5.630939170309886)"
